$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '21.292.01'
Set-TextValue $ws.Range("E2") '  +4.49%  '

Set-TextValue $ws.Range("D3") '1.544.85'
Set-TextValue $ws.Range("E3") '  +5.78%  '

Set-TextValue $ws.Range("D4") '1.004'
Set-TextValue $ws.Range("E4") '  -0.56%  '

Set-TextValue $ws.Range("D5") '0.9585'

Set-TextValue $ws.Range("D6") '282.57'
Set-TextValue $ws.Range("E6") '  +2.88%  '

Set-TextValue $ws.Range("D7") '0.3635'
Set-TextValue $ws.Range("E7") '  -0.32%  '

Set-TextValue $ws.Range("D8") '0.3195'
Set-TextValue $ws.Range("E8") '  +4.13%  '

Set-TextValue $ws.Range("D9") '40.98'
Set-TextValue $ws.Range("E9") '  +3.35%  '

Set-TextValue $ws.Range("D10") '1.097'
Set-TextValue $ws.Range("E10") '  +6.18%  '

Set-TextValue $ws.Range("D11") '0.06828'
Set-TextValue $ws.Range("E11") '  +3.90%  '

Set-TextValue $ws.Range("D12") '0.9982'
Set-TextValue $ws.Range("E12") '  -0.31%  '

Set-TextValue $ws.Range("D13") '5.686'
Set-TextValue $ws.Range("E13") '  +4.94%  '

Set-TextValue $ws.Range("D14") '18.82'
Set-TextValue $ws.Range("E14") '  +5.55%  '

Set-TextValue $ws.Range("D15") '6.367'
Set-TextValue $ws.Range("E15") '  +3.96%  '

Set-TextValue $ws.Range("D16") '0.00001052'
Set-TextValue $ws.Range("E16") '  +2.85%  '

Set-TextValue $ws.Range("B17") 'WrappedEther'
Set-TextValue $ws.Range("C17") 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D17") '1.542.94'
Set-TextValue $ws.Range("E17") '  +5.63%  '

Set-TextValue $ws.Range("B18") 'Dai'
Set-TextValue $ws.Range("C18") 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D18") '0.9592'
Set-TextValue $ws.Range("E18") '  -0.90%  '

Set-TextValue $ws.Range("D19") '0.06065'
Set-TextValue $ws.Range("E19") '  +4.74%  '

Set-TextValue $ws.Range("D20") '72.55'
Set-TextValue $ws.Range("E20") '  +4.35%  '

Set-TextValue $ws.Range("D21") '5.708'
Set-TextValue $ws.Range("E21") '  +5.08%  '

Set-TextValue $ws.Range("D22") '15.07'
Set-TextValue $ws.Range("E22") '  +4.55%  '

Set-TextValue $ws.Range("E23") '  +4.65%  '

Set-TextValue $ws.Range("D24") '2.313'
Set-TextValue $ws.Range("E24") '  +2.88%  '

Set-TextValue $ws.Range("D25") '21.353.15'
Set-TextValue $ws.Range("E25") '  +4.59%  '

Set-TextValue $ws.Range("D26") '148.29'
Set-TextValue $ws.Range("E26") '  +4.81%  '

Set-TextValue $ws.Range("D27") '2.223'
Set-TextValue $ws.Range("E27") '  +6.98%  '

Set-TextValue $ws.Range("D28") '17.71'
Set-TextValue $ws.Range("E28") '  +3.52%  '

Set-TextValue $ws.Range("D29") '1.702.47'
Set-TextValue $ws.Range("E29") '  +5.35%  '

Set-TextValue $ws.Range("D30") '117.97'
Set-TextValue $ws.Range("E30") '  +5.36%  '

Set-TextValue $ws.Range("D31") '4.029'
Set-TextValue $ws.Range("E31") '  +5.48%  '

Set-TextValue $ws.Range("D32") '5.250'
Set-TextValue $ws.Range("E32") '  +7.52%  '

Set-TextValue $ws.Range("D33") '0.8546'
Set-TextValue $ws.Range("E33") '  +8.37%  '

Set-TextValue $ws.Range("D34") '0.08017'
Set-TextValue $ws.Range("E34") '  +1.77%  '

Set-TextValue $ws.Range("D35") '1.507'
Set-TextValue $ws.Range("E35") '  -1.18%  '

Set-TextValue $ws.Range("D36") '4.986'
Set-TextValue $ws.Range("E36") '  +6.79%  '

Set-TextValue $ws.Range("E37") '  +6.13%  '

Set-TextValue $ws.Range("D38") '0.05888'
Set-TextValue $ws.Range("E38") '  +3.28%  '

Set-TextValue $ws.Range("D39") '0.02109'
Set-TextValue $ws.Range("E39") '  +4.28%  '

Set-TextValue $ws.Range("D40") '10.75'
Set-TextValue $ws.Range("E40") '  +4.32%  '

Set-TextValue $ws.Range("D41") '7.753'
Set-TextValue $ws.Range("E41") '  +3.80%  '

Set-TextValue $ws.Range("D42") '0.1920'
Set-TextValue $ws.Range("E42") '  +3.57%  '

Set-TextValue $ws.Range("D43") '0.9596'
Set-TextValue $ws.Range("E43") '  +0.25%  '

Set-TextValue $ws.Range("D44") '0.5474'
Set-TextValue $ws.Range("E44") '  +4.26%  '

Set-TextValue $ws.Range("D45") '12.53'
Set-TextValue $ws.Range("E45") '  +5.80%  '

Set-TextValue $ws.Range("D46") '3.581'
Set-TextValue $ws.Range("E46") '  +2.74%  '

Set-TextValue $ws.Range("D47") '0.5471'
Set-TextValue $ws.Range("E47") '  +7.00%  '

Set-TextValue $ws.Range("D48") '121.96'
Set-TextValue $ws.Range("E48") '  +4.29%  '

Set-TextValue $ws.Range("D49") '1.877'
Set-TextValue $ws.Range("E49") '  +7.46%  '

Set-TextValue $ws.Range("D50") '0.06635'
Set-TextValue $ws.Range("E50") '  +3.47%  '

Set-TextValue $ws.Range("D51") '70.13'
Set-TextValue $ws.Range("E51") '  +6.08%  '
